$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 12-14 do not exist yet as styled cells (unlike rows 9-11 which already
# carry the table's body formatting even while empty). Copy the formatting
# from the existing body row 9 onto the new rows before writing any values,
# so the new rows pick up the same look (bold centered B/C columns, left
# aligned wrapped D column, centered F column) as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("B9:D9").Copy() | Out-Null
$ws.Range("B12:D14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F9").Copy() | Out-Null
$ws.Range("F11:F14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 9: existing row - description/header text swap places, row shrinks,
#            and the previously-empty test-case count is filled in. ---------
$ws.Range("B9").Value = "(TS_006)`nHeader_Contact Support"
$ws.Range("C9").Value = "FRS"
$ws.Range("D9").Value = "Validte the functionality of the Dropdown Header > Contact Support."
$ws.Range("F9").Value = 3
$ws.Rows.Item(9).RowHeight = 30

# --- Row 10: used to be a blank spacer/subtotal row, now a data row --------
$ws.Range("B10").Value = "(TS_007)`nHeader_Contact Support"
$ws.Range("C10").Value = "FRS"
$ws.Range("D10").Value = "Validte the functionality of the Dropdown Header > Logout."
$ws.Range("F10").Value = 1
$ws.Rows.Item(10).RowHeight = 30

# --- Row 11: used to be a blank spacer row, now a data row -----------------
$ws.Range("B11").Value = "(TS_008)`nHome Page_ Dash Board"
$ws.Range("C11").Value = "FRS"
$ws.Range("D11").Value = "Validte the functionality of the Home Page > Dash Board."
$ws.Range("F11").Value = 34
$ws.Rows.Item(11).RowHeight = 30

# --- Row 12: new data row ---------------------------------------------------
$ws.Range("B12").Value = "(TS_009)`nHome Page_ Recent Activity"
$ws.Range("C12").Value = "FRS"
$ws.Range("D12").Value = "Validte the functionality of the Home Page > Recent Activity."
$ws.Range("F12").Value = 2
$ws.Rows.Item(12).RowHeight = 45

# --- Row 13: new data row ---------------------------------------------------
$ws.Range("B13").Value = "(TS_010)`nHome Page_ All Expenses"
$ws.Range("C13").Value = "FRS"
$ws.Range("D13").Value = "Validte the functionality of the Home Page > All Expenses."
$ws.Range("F13").Value = 15
$ws.Rows.Item(13).RowHeight = 45

# --- Row 14: new data row ---------------------------------------------------
$ws.Range("B14").Value = "(TS_011)`nHome Page_ Groups"
$ws.Range("C14").Value = "FRS"
$ws.Range("D14").Value = "Validte the functionality of the Home Page > Groups."
$ws.Range("F14").Value = 21
$ws.Rows.Item(14).RowHeight = 30

# --- Totals block (rows 15-17) ---------------------------------------------
$ws.Range("F15").Formula = "=SUM(F4:F14)"
$ws.Range("F16").Value = 300
$ws.Range("F17").Formula = "=F16-F15"

# --- View state: active cell / scroll position ------------------------------
$ws.Range("F15").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
